$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 45455770
$ws.Range("I18").Value = 55556740
$ws.Range("J18").Value = 1385
$ws.Range("K18").Value = 55556740
$ws.Range("L18").Value = 1385
$ws.Range("M18").Value = -55556456
$ws.Range("N18").Value = -1953

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1947

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2099.9285
$ws.Range("I98").Value = 1769.3
$ws.Range("K98").Value = 1769.3
$ws.Range("M98").Value = -271.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1834.2858
$ws.Range("I100").Value = 1660
$ws.Range("J100").Value = 1888.75
$ws.Range("K100").Value = 1660
$ws.Range("L100").Value = 1888.75
$ws.Range("M100").Value = -1119
$ws.Range("N100").Value = -2970.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2099.9285
$ws.Range("I122").Value = 1769.3
$ws.Range("K122").Value = 5307.9
$ws.Range("M122").Value = -2857.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 14927053
$ws.Range("I138").Value = 23810802
$ws.Range("J138").Value = 2355.44
$ws.Range("K138").Value = 71432406
$ws.Range("L138").Value = 7066.32
$ws.Range("M138").Value = -71427266
$ws.Range("N138").Value = -17346.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 293.53
$ws.Range("I32").Value = 269.08334
$ws.Range("J32").Value = 880.25
$ws.Range("K32").Value = 269.08334
$ws.Range("L32").Value = 880.25
$ws.Range("M32").Value = 17.91665999999998
$ws.Range("N32").Value = -1454.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7222.067
$ws.Range("I61").Value = 7585.926
$ws.Range("K61").Value = 7585.926
$ws.Range("M61").Value = -7373.926

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3472.2727
$ws.Range("I63").Value = 3149.375
$ws.Range("J63").Value = 4333.3335
$ws.Range("K63").Value = 3149.375
$ws.Range("L63").Value = 4333.3335
$ws.Range("M63").Value = -2463.375
$ws.Range("N63").Value = -5705.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3472.2727
$ws.Range("I66").Value = 3149.375
$ws.Range("J66").Value = 4333.3335
$ws.Range("K66").Value = 15746.875
$ws.Range("L66").Value = 21666.6675
$ws.Range("M66").Value = -12314.875
$ws.Range("N66").Value = -28530.6675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7222.067
$ws.Range("I136").Value = 7585.926
$ws.Range("K136").Value = 22757.778
$ws.Range("M136").Value = -20207.778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 41918.668
$ws.Range("I22").Value = 226
$ws.Range("J22").Value = 166996.67
$ws.Range("K22").Value = 226
$ws.Range("L22").Value = 166996.67
$ws.Range("M22").Value = -53
$ws.Range("N22").Value = -167342.67

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2329.875
$ws.Range("I107").Value = 2022.7693
$ws.Range("J107").Value = 3660.6667
$ws.Range("K107").Value = 2022.7693
$ws.Range("L107").Value = 3660.6667
$ws.Range("M107").Value = -102.7692999999999
$ws.Range("N107").Value = -7500.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2601.9546
$ws.Range("I31").Value = 1972.5714
$ws.Range("J31").Value = 3703.375
$ws.Range("K31").Value = 1972.5714
$ws.Range("L31").Value = 3703.375
$ws.Range("M31").Value = -1677.5714
$ws.Range("N31").Value = -4293.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2601.9546
$ws.Range("I34").Value = 1972.5714
$ws.Range("J34").Value = 3703.375
$ws.Range("K34").Value = 1972.5714
$ws.Range("L34").Value = 3703.375
$ws.Range("M34").Value = -1770.5714
$ws.Range("N34").Value = -4107.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6164.1963
$ws.Range("I58").Value = 5798.4326
$ws.Range("J58").Value = 7130.857
$ws.Range("K58").Value = 5798.4326
$ws.Range("L58").Value = 7130.857
$ws.Range("M58").Value = -5595.4326
$ws.Range("N58").Value = -7536.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 49593.332
$ws.Range("J109").Value = 49593.332
$ws.Range("L109").Value = 49593.332
$ws.Range("N109").Value = -51673.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 8504.348
$ws.Range("I134").Value = 9684.4375
$ws.Range("K134").Value = 29053.3125
$ws.Range("M134").Value = -26518.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6164.1963
$ws.Range("I136").Value = 5798.4326
$ws.Range("J136").Value = 7130.857
$ws.Range("K136").Value = 17395.2978
$ws.Range("L136").Value = 21392.571
$ws.Range("M136").Value = -14845.2978
$ws.Range("N136").Value = -26492.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 76718
$ws.Range("I137").Value = 55000
$ws.Range("J137").Value = 87577
$ws.Range("K137").Value = 55000
$ws.Range("L137").Value = 87577
$ws.Range("M137").Value = -49900
$ws.Range("N137").Value = -97777

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 68755.44
$ws.Range("J141").Value = 74249.64
$ws.Range("L141").Value = 74249.64
$ws.Range("N141").Value = -84609.64

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 456.47058
$ws.Range("I5").Value = 397.3793
$ws.Range("J5").Value = 799.2
$ws.Range("K5").Value = 1192.1379
$ws.Range("L5").Value = 2397.6
$ws.Range("M5").Value = -1080.1379
$ws.Range("N5").Value = -2621.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 12836.75
$ws.Range("I82").Value = 8216.5
$ws.Range("K82").Value = 24649.5
$ws.Range("M82").Value = -24243.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 12836.75
$ws.Range("I85").Value = 8216.5
$ws.Range("K85").Value = 24649.5
$ws.Range("M85").Value = -23245.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 393.8889
$ws.Range("J97").Value = 346.5
$ws.Range("L97").Value = 1039.5
$ws.Range("N97").Value = -2031.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 5089.8335
$ws.Range("I132").Value = 4055
$ws.Range("K132").Value = 36495
$ws.Range("M132").Value = -33965

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 456.47058
$ws.Range("I135").Value = 397.3793
$ws.Range("J135").Value = 799.2
$ws.Range("K135").Value = 3576.4137
$ws.Range("L135").Value = 7192.8
$ws.Range("M135").Value = -1041.4137
$ws.Range("N135").Value = -12262.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2028.3572
$ws.Range("I97").Value = 2090.3157
$ws.Range("J97").Value = 1897.5555
$ws.Range("K97").Value = 2090.3157
$ws.Range("L97").Value = 1897.5555
$ws.Range("M97").Value = -1594.3157
$ws.Range("N97").Value = -2889.5555

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5003.5884
$ws.Range("I122").Value = 4575.9287
$ws.Range("K122").Value = 13727.7861
$ws.Range("M122").Value = -11277.7861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4560.1577
$ws.Range("I40").Value = 4637.16
$ws.Range("J40").Value = 4412.077
$ws.Range("K40").Value = 4637.16
$ws.Range("L40").Value = 4412.077
$ws.Range("M40").Value = -4501.16
$ws.Range("N40").Value = -4684.077

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 938.7857
$ws.Range("I93").Value = 1033.909
$ws.Range("J93").Value = 590
$ws.Range("K93").Value = 1033.909
$ws.Range("L93").Value = 590
$ws.Range("M93").Value = 214.0909999999999
$ws.Range("N93").Value = -3086

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2355.9092
$ws.Range("I100").Value = 2110.5
$ws.Range("J100").Value = 2650.4
$ws.Range("K100").Value = 2110.5
$ws.Range("L100").Value = 2650.4
$ws.Range("M100").Value = -1569.5
$ws.Range("N100").Value = -3732.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 117664.664
$ws.Range("J139").Value = 117664.664
$ws.Range("L139").Value = 117664.664
$ws.Range("N139").Value = -127944.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 9538.625
$ws.Range("I126").Value = 5465
$ws.Range("K126").Value = 16395
$ws.Range("M126").Value = -13925

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3282.1428
$ws.Range("I132").Value = 2201.3076
$ws.Range("K132").Value = 6603.9228
$ws.Range("M132").Value = -4073.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2008.8975
$ws.Range("I136").Value = 1915.1351
$ws.Range("K136").Value = 5745.4053
$ws.Range("M136").Value = -3195.4053
